$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(202, 'Objęto go ochroną w celach naukowo-dydaktycznych, ze względu na atrakcyjny pokrój drzewa oraz budowę morfologiczną liści i kory drzewa, a także ze względu na pokaźne rozmiary.', 'polish')
    ,@(203, 'Stan zdrowotny drzewa jest zadowalający, aczkolwiek nie jest ono zachowane w całości – pozostała głównie kłoda i ok. 40% masy korony, brak jest partii wierzchołkowej, spowodowany złamaniem przewodnika oraz wykonaniem cięć redukcyjnych.', 'polish')
    ,@(204, 'Drużynowy konkurs skoków narciarskich na normalnej skoczni powrócił po tym, jak nie odbył się on na Mistrzostwach Świata w 2003.', 'polish')
    ,@(205, 'Miał obszerną bibliotekę specjalistycznych książek, literatury pięknej, a także dzieł słowiańskich.', 'polish')
    ,@(206, 'W czasie wojny domowej miejscowość znalazła się w rejonie działania antyrządowych bojówek.', 'polish')
    ,@(207, 'W 2020 roku miejscowe władze odkryły pokaźny arsenał terrorystów, w tym karabiny automatyczne, snajperskie, pociski moździerzowe, rakiety, urządzenia wybuchowe i znaczną ilość amunicji', 'polish')
    ,@(208, 'Zgodnie z tradycją co roku w mieście organizowany jest festiwal ziemniaków.', 'polish')
    ,@(209, 'Astronautyka obejmuje poznawanie oraz analizowanie warunków i zjawisk towarzyszących lotom statków kosmicznych. ', 'polish')
    ,@(210, 'Odkrywa możliwości techniczne ich realizacji oraz bada oddziaływanie warunków lotu na psychofizyczny stan astronautów i możliwość ich adaptacji do tych warunków.', 'polish')
    ,@(211, 'W rozwiązywaniu zagadnień z zakresu astronautyki korzysta się z osiągnięć wielu gałęzi nauki (aerodynamiki, magnetohydrodynamiki, aeronomii, astrofizyki, planetologii, techniki rakietowej i innych).', 'polish')
    ,@(212, 'Idea podróży kosmicznych nurtowała ludzkie umysły od stuleci, ale pozostawała marzeniem aż do czasu zbudowania potężnych rakiet zdolnych unieść użyteczny ładunek daleko w przestrzeń.', 'polish')
    ,@(213, 'W ciągu niespełna trzydziestu lat sondy kosmiczne przeleciały obok wszystkich planet Układu Słonecznego.', 'polish')
    ,@(214, 'Tańce mogą być podzielone i opisane ze względu na rodzaj choreografii, rodzaj ruchów lub historyczne pochodzenie.', 'polish')
    ,@(215, 'Przed wynalezieniem pisma taniec był ważnym czynnikiem międzypokoleniowego przekazu historycznego.', 'polish')
    ,@(216, 'Najprostszy podział, jakiego można dokonać, to wyodrębnienie trzech odmiennych sfer działalności ludzkiej, w których się rozwinął: religia, rozrywka i sztuka.', 'polish')
    ,@(217, 'Poczta w znaczeniu nowożytnym, będąca częścią europejskiej sieci pocztowej, powstała za czasów króla Zygmunta Augusta.', 'polish')
    ,@(218, 'Zaprezentowany dwa lata później Szantaż był pierwszym brytyjskim filmem dźwiękowym.', 'polish')
    ,@(219, 'Arboretum w Rogowie – należy do grupy najcenniejszych, najbogatszych w gatunki i odmiany drzew i krzewów tego typu ogrodów w Europie.', 'polish')
    ,@(220, 'Obok jeziora Żółkin przebiega XII trasa rowerowa Bełchatowskiego Stowarzyszenia Załogi Rowerowej "Zgrzyt".', 'polish')
    ,@(221, 'Nadleśnictwo położone jest w północno-zachodniej części terenu województwa wielkopolskiego na terenie 5 powiatów: obornickim, szamotulskim, czarnkowsko-trzcianeckim, chodzieskim i poznańskim.', 'polish')
    ,@(222, 'Macedonia was an ancient kingdom on the periphery of Archaic and Classical Greece, and later the dominant state of Hellenistic Greece.', 'english')
    ,@(223, 'Canada bans the sale of assault-style firearms after a mass killing in Nova Scotia leaves twenty-two victims dead.', 'english')
    ,@(224, 'The remoteness of the Glensanda settlement is such that there are no road, rail, or marked footway links across the granite mountain, moor, heather and peat bog of the private Glensanda estate.', 'english')
    ,@(225, 'The only practical access is by boat from the shores of Loch Linnhe.', 'english')
    ,@(226, ' The band was playing the club circuit in London, and he was invited to join them for a few songs.', 'english')
    ,@(227, 'There are roads and private residences along much of the lake, except the southeastern side where the Caribou Bog complex borders the lake.', 'english')
    ,@(228, 'The squadron returned to the United States in December 1945 and was inactivated in March 1946, and its personnel and equipment transferred to another organization.', 'english')
    ,@(229, 'Organdy or organdie is the sheerest and crispest cotton cloth made.', 'english')
    ,@(230, 'The latter two finishes are more popular for summer wear and draped apparel whereas the first is more popular for loose apparel and home textiles such as dresses and curtains.', 'english')
    ,@(231, 'The film also denounces gang violence and presents meaningful solutions from former gang-members to stop this problem.', 'english')
    ,@(232, 'Her collected writings and recollections of that period offer a rare first hand English language accounts of life in wartime Japan.', 'english')
    ,@(233, 'Zanskar, ringed by high Himalayan mountains in northwest India, one of the most remote places on the planet, has been safe until now.', 'english')
    ,@(234, 'While they complete the school they are also placing local children in other schools and monasteries in the city of Manali and beyond. ', 'english')
    ,@(235, 'After this, the Serbian player pushed Shengelia and the Georgian forward reacted with a punch. ', 'english')
    ,@(236, 'After this, all the players who were in the bench came into the court and started the brawl.', 'english')
    ,@(237, 'Waller plans on sponsoring legislation to require disclosure on clothes collection bins operated by for-profit companies, to include a new derivative of methamphetamine to Colorado''s list of banned drugs, and to exclude dry-ice bombs from a felony list of explosives, and to allow homeschooled students to enroll in college classes.', 'english')
    ,@(238, 'On his dismissal as governor, he said of the State treasury that he left it empty because he met it empty.', 'english')
    ,@(239, 'In June 2014, its assets and infrastructure were purchased by a new company.', 'english')
    ,@(240, 'The village has a population of 589.', 'english')
    ,@(241, 'That is, the information is intended to be quickly found when needed.', 'english')
    ,@(242, 'La biología se ocupa tanto de la descripción de las características y los comportamientos de los organismos individuales, como de las especies en su conjunto, así como de la reproducción de los seres vivos y de las interacciones entre ellos y el entorno.', 'spanish')
    ,@(243, 'La biología moderna se divide en sub-disciplinas según los tipos de organismos y la escala en que se los estudia. ', 'spanish')
    ,@(244, 'El virus se transmite generalmente de persona a persona a través de las pequeñas gotas de saliva.', 'spanish')
    ,@(245, 'Se han cerrado colegios y universidades en más de 124 países, lo que ha afectado a más de 2200 millones de estudiantes.', 'spanish')
    ,@(246, 'La mayoría de las provincias implementaron cierres de escuelas y guarderías, prohibiciones de grandes reuniones, así como el cierre de varios lugares de ocio y entretenimiento a mediados de marzo.', 'spanish')
    ,@(247, 'En Corea del Norte no se ha registrado ningún caso oficial, lo cual fue muy cuestionado debido a sus fronteras con China y Corea del Sur.', 'spanish')
    ,@(248, 'Puedes avisarnos de un artículo de actualidad que falte.', 'spanish')
    ,@(249, 'Los componentes del movimiento, aunque conservan su independencia en los límites del estatuto del movimiento, actúan siempre de conformidad con sus principios fundamentales y colaboran entre sí en el desempeño de sus tareas respectivas y para realizar su misión común.', 'spanish')
    ,@(250, 'En la segunda mitad del siglo XIX, la rápida evolución de la tecnología de las armas de fuego ocasionó un dramático incremento del número de muertos y de heridos en tiempo de guerra.', 'spanish')
    ,@(251, 'El símbolo debía ser sencillo, identificable a distancia, conocido por todos e idéntico para amigos y adversarios. ', 'spanish')
    ,@(252, 'En el primer Convenio de Ginebra, se eligió la cruz roja sobre fondo blanco como emblema distintivo único.  ', 'spanish')
    ,@(253, ' Considera también revisar los espacios de nombres distintos al principal, especialmente los anexos, aquí.', 'spanish')
    ,@(254, 'En la ciénega se reúnen varios arroyos con los que dan vida al río Lauca.', 'spanish')
    ,@(255, 'Una vez calculados para cada mes, ambos valores son calculados para todo el año y pueden ser leídos en la columna vertical al lado derecho del diagrama.', 'spanish')
    ,@(256, 'Sus contribuciones en econometría incluir el estudio de riesgo y seguros en los países en desarrollo.', 'spanish')
    ,@(257, 'Tepito es una estación del Metro de Ciudad de México perteneciente a la línea B, localizada en el barrio de Tepito.', 'spanish')
    ,@(258, 'Causa pudrición de las raíces en muchas especies de plantas.', 'spanish')
    ,@(259, 'Carne blanquecina y firme en el sombrero, leñosa y fibrosa en el pie, de sabor suave en los ejemplares jóvenes, amargo y desagradable en los ejemplares adultos, y de olor fuerte no muy agradable.  ', 'spanish')
    ,@(260, 'Por lo tanto los dos sustratos de esta enzima son nitrito, y iones hidrógeno; mientras que sus tres productos son óxido nítrico, nitrato, y agua.', 'spanish')
    ,@(261, 'El anfitrión provee nutrientes orgánicos al hongo, y la actividad metabólica dentro del complejo es considerablemente mayor que afuera.', 'spanish')
    ,@(262, 'On emploie aussi le terme de suçoir intercellulaire (intracellulaire en apparence), pour désigner l''expansion développée par ces parasites et qui traverse la paroi pectocellulosique en repoussant légèrement la membrane de la cellule, sans la rompre. ', 'french')
    ,@(263, 'L''œuvre poétique latine d''Arrius Nurus est une des plus importantes et des plus sensibles du vingtième siècle.', 'french')
    ,@(264, 'Alors qu''il joue en minime deuxième année, Dugarry rejoint les Girondins de Bordeaux et débute à 17 ans en équipe première en championnat de France.', 'french')
    ,@(265, 'Parmi les œuvres du compositeur qui ont survécu, il faut citer une pièce en do majeur qui pourrait être la plus ancienne suite pour trois violons de l''histoire de la musique, une pièce virtuose requérant une maîtrise parfaite de la technique instrumentale.  ', 'french')
    ,@(266, 'Les travaux de construction de la nouvelle église vont débuter début 2020.  ', 'french')
    ,@(267, 'Les procès en trahison étaient une extension des procès pour sédition de 1792 et 1793 contre les parlementaires réformateurs d''Angleterre et d''Écosse.', 'french')
    ,@(268, 'Cette position indispose le Ministère et son propre parti.', 'french')
    ,@(269, 'L''année suivante, l''écurie belge n''est que suppléant sur la liste des engagés des 24 Heures du Mans', 'french')
    ,@(270, 'On dit que deux droites sont sécantes si elles ont un unique point commun.  ', 'french')
    ,@(271, 'Comme tous les aliments en saumure, sa saveur est forte, très acide, ce qui le rend idéal pour sa consommation comme apéritif accompagné de vin.  ', 'french')
    ,@(272, 'Les données de stock sont partagées et mises à jour régulièrement.', 'french')
    ,@(273, 'Les deux parties disposent ainsi des mêmes informations pour leur planification.', 'french')
    ,@(274, 'Le fournisseur optimise le stock en anticipant la production; gérant les dates et quantités de livraison pour chaque pièce.', 'french')
    ,@(275, 'Amélioration de la fiabilité des livraisons et optimisation de la sécurité d’approvisionnement. ', 'french')
    ,@(276, 'Réduction des étapes de planification et simplification du travail des gestionnaires. ', 'french')
    ,@(277, 'La chasse illégale pourrait néanmoins constituer une menace non négligeable.', 'french')
    ,@(278, 'Lors du recensement de 2010, la municipalité comptait 131 719 habitan.', 'french')
    ,@(279, 'Dans cette discipline, elle a aussi été deux fois quart de finaliste en Grand Chelem.', 'french')
    ,@(280, 'Chaque pays a été autorisé à participer avec deux équipes mixte de quatre participants, composée de deux femmes et de deux hommes', 'french')
    ,@(281, 'Une psyché est un grand miroir mobile monté sur un châssis que l''on peut incliner à volonté pour se regarder en pied.̟  ', 'french')
    ,@(282, 'Das Möbelstück Psyche wird im Damen Conversations Lexikon von 1838 als Ankleidespiegel bezeichnet, der sich im Ankleidezimmer befand.', 'german')
    ,@(283, 'Das Hotelgebäude wurde ursprünglich als Apartmenthaus Ende des 18.', 'german')
    ,@(284, 'Dabei sind die Bestandteile alter Kurs, Manövereinleitung, schnelles Schiften, Manöverende und neuer Kurs klar zu identifizieren.', 'german')
    ,@(285, 'Die scharfe Abgrenzung von Schifte und Halse ist eine Besonderheit der deutschen Segelausbildung.', 'german')
    ,@(286, 'Der Künstler wohnt und arbeitet in Göteborg.', 'german')
    ,@(287, 'Alberto Barberis geriet in Vergessenheit und verstarb im Jahr 1976 im Alter von 93 Jahren in einem Altersheim.  ', 'german')
    ,@(288, 'Das Mineral ist damit chemisch gesehen ein Zink-Aluminium-Oxid oder auch Zinkaluminat.', 'german')
    ,@(289, 'Das Mineral ist durchscheinend bis undurchsichtig und zeigt auf den Oberflächen der blaugrünen, gelbbraunen oder grünlichschwarzen Kristalle einen fettigen Glasglanz.  ', 'german')
    ,@(290, 'Als Erz ist Gahnit ohne wirtschaftliche Bedeutung.', 'german')
    ,@(291, 'Gelegentlich wird er aber von versierten Sammlern in Facettenform geschliffen.', 'german')
    ,@(292, 'Über die Merwede verlaufen flussabwärts von Ost nach West die folgenden Brücken.', 'german')
    ,@(293, 'Dazu zeigt man den Fortsetzungssatz mit Hilfe des Interpolationssatzes zunächst für stetige Funktionen, die eine abgeschlossene Teilmenge eines normalen topologischen Raumes in das Intervall {\displaystyle [0,1]}[0,1] abbilden.', 'german')
    ,@(294, 'Wie im gesamten deutschen Sprachraum ist auch in Baden der schwarze Talar mit weißem Beffchen die Amtstracht der evangelischen Geistlichen.', 'german')
    ,@(295, 'Im Areal dieser Siedlungen ist als Relikt des Leitenwaldes ein parkartiges Waldstück unverbaut geblieben, das Napoleonwald genannt wird.  ', 'german')
    ,@(296, 'Erste Einträge über die Langenthaler Fasnacht sind ab 1864 bekannt.  ', 'german')
    ,@(297, 'Hüpfschritt zu einem speziellen Rhythmus in der Gruppe. ', 'german')
    ,@(298, 'Die Friedhofskapelle ist heute eine Aufbahrungshalle und in ihrer ursprünglichen Bauart und Größe erhalten.', 'german')
    ,@(299, 'Der Warzen-Teppichhai ist ein mittelgroßer, bodenlebender Hai mit einer durchschnittlichen Länge von 0,65 m, die Maximalgröße reicht bis etwa 0,90 m.', 'german')
    ,@(300, 'Verwandte bicyclische Verbindungen sind Norbornen, welches dasselbe Kohlenstoffgerüst, aber eine Doppelbindung besitzt, und das Norbornadien, ebenfalls mit demselben Gerüst, aber zwei Doppelbindungen und somit ein Dien.', 'german')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

# Match the author's final column-A width resize and cursor position
$ws.Columns.Item(1).ColumnWidth = 223.45
$null = $ws.Range("B301").Select()
